# Actualización automática 2025-10-08 14:30:09
# Insert a new advisor/client row ("PAUTA ASTUDILLO JULIO HERNAN") above the
# existing "VIEJO RIVAS MAYRA ANABELLE" row on both data sheets, pushing the
# totals row down and bumping the "X de N" denominators accordingly.

$wb = $excel.ActiveWorkbook

function Insert-PautaRow {
    param($SheetName, $LastCol, $HasTotalsText, $HTotalsCol)

    $ws = $wb.Worksheets.Item($SheetName)

    # Row 13 currently holds "VIEJO RIVAS MAYRA ANABELLE"; insert a fresh row
    # above it so that row becomes row 14, and the old totals row (14) becomes 15.
    $ws.Rows("13:13").Insert()

    # Column A keeps repeating the same advisor name for every client row.
    $ws.Cells.Item(13, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
    $ws.Cells.Item(13, 2).Value = "PAUTA ASTUDILLO JULIO HERNAN"

    for ($c = 3; $c -le $LastCol; $c++) {
        $ws.Cells.Item(13, $c).Value = 0
    }

    if ($HasTotalsText) {
        # Row 15 is the former "0 de 12" totals row; bump the denominator to 13
        # now that there is one more client row feeding it.
        for ($c = 3; $c -le $LastCol; $c++) {
            if ($c -eq $HTotalsCol) {
                $ws.Cells.Item(15, $c).Value = "1 de 13"
            } else {
                $ws.Cells.Item(15, $c).Value = "0 de 13"
            }
        }
    }
}

Insert-PautaRow "VENTAS POR GRUPO" 18 $true 8
Insert-PautaRow "VENTA MENSUAL" 7 $false 0
